$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 102 ---
# A102: date style (matches existing style used by A100, numFmtId 14 / border 4)
$ws.Range("A100").Copy()
$ws.Range("A102").PasteSpecial(-4122)
$ws.Range("A102").Value = 45895

# B102 / C102: time style (matches existing style used by B100/C100, numFmtId 20 / border 4)
$ws.Range("B100").Copy()
$ws.Range("B102").PasteSpecial(-4122)
$ws.Range("B102").Value = 0.5625

$ws.Range("C100").Copy()
$ws.Range("C102").PasteSpecial(-4122)
$ws.Range("C102").Value = 0.625

# D102: new style - numFmtId 16 ("d-mmm") with the same thin border as the rest of the row
$ws.Range("D102").Borders.LineStyle = 0
$ws.Range("D102").Borders.LineStyle = 1
$ws.Range("D102").Borders.Weight = 2
$ws.Range("D102").NumberFormat = "d-mmm"

# --- Row 103 ---
$ws.Range("A100").Copy()
$ws.Range("A103").PasteSpecial(-4122)
$ws.Range("A103").Value = 45897

$ws.Range("B100").Copy()
$ws.Range("B103").PasteSpecial(-4122)
$ws.Range("B103").Value = 0.58333333333333337

$ws.Range("C100").Copy()
$ws.Range("C103").PasteSpecial(-4122)
$ws.Range("C103").Value = 0.625

$excel.CutCopyMode = 0

$ws.Range("E103").Select()
